$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" banner timestamp (A1): 20:22 -> 20:52
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 20:52"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 489646
$ws.Range("C4").Value = 21080
$ws.Range("D4").Value = 26777
$ws.Range("E4").Value = 444835
$ws.Range("G4").Value = 1343
$ws.Range("H4").Value = 18034

# Row 14: Suiza - refreshed totals
$ws.Range("D14").Value = 11100
$ws.Range("E14").Value = 12449
$ws.Range("G14").Value = 54
$ws.Range("H14").Value = 1002

# Row 17: Brasil - refreshed totals
$ws.Range("B17").Value = 19638
$ws.Range("C17").Value = 1493
$ws.Range("E17").Value = 18408
$ws.Range("G17").Value = 103
$ws.Range("H17").Value = 1057

# Rows 57-58: Egipto's case count overtakes Argelia's, so Egipto now sorts
# above Argelia in the (descending, by total cases) ranking. Egipto gets the
# refreshed numbers at row 57; Argelia's (unchanged) numbers drop to row 58.
$ws.Range("A57").Value = "Egipto"
$ws.Range("B57").Value = 1794
$ws.Range("C57").Value = 95
$ws.Range("D57").Value = 384
$ws.Range("E57").Value = 1275
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 17
$ws.Range("H57").Value = 135

$ws.Range("A58").Value = "Argelia"
$ws.Range("B58").Value = 1761
$ws.Range("C58").Value = 95
$ws.Range("D58").Value = 405
$ws.Range("E58").Value = 1100
$ws.Range("F58").Value = 46
$ws.Range("G58").Value = 21
$ws.Range("H58").Value = 256

# Row 59: Islandia - refreshed totals
$ws.Range("E59").Value = 917
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 7

# Row 139: Barbados - refreshed totals
$ws.Range("B139").Value = 67
$ws.Range("C139").Value = 1
$ws.Range("E139").Value = 52
